$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.138055562973022
$ws.Range("B1").Value = 2.585058689117432
$ws.Range("C1").Value = 6.919022083282471
$ws.Range("D1").Value = 2.091269731521606
$ws.Range("E1").Value = 1.235441565513611
